{"js": "// Make the label at the start of each \"Label: value\" paragraph bold, and\n// normalize the label to Title Case (matching the author's commit: \"modify\n// the paragraphs to make the title bold on each section\").\n//\n// Each target paragraph currently holds a single run \"Label: value\". We\n// split it into two runs: a bold run with the (title-cased) label text, and\n// a plain run with \": value\" (colon, space, and the rest unchanged).\n\n// Canonical label text (title-cased), keyed by the lower-cased original\n// label so the match is resilient to minor case differences.\nconst LABELS = [\n  \"Project title\",\n  \"Project description\",\n  \"Business driver\",\n  \"Business value\",\n  \"Business risk\",\n  \"Budget expense\",\n  \"Internal hours\",\n  \"External hours\",\n  \"Solutions involvement\",\n  \"Solutions hours\",\n  \"PMO involvement\",\n  \"PMO hours\",\n  \"Total expected hours\",\n];\n\nfunction titleCase(label) {\n  return label\n    .split(\" \")\n    .map((w) => (w.toUpperCase() === \"PMO\" ? \"PMO\" : w.charAt(0).toUpperCase() + w.slice(1).toLowerCase()))\n    .join(\" \");\n}\n\nconst LABEL_MAP = new Map(LABELS.map((l) => [l.toLowerCase(), titleCase(l)]));\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const text = para.text;\n  const sep = text.indexOf(\": \");\n  if (sep === -1) continue;\n\n  const rawLabel = text.slice(0, sep);\n  const newLabel = LABEL_MAP.get(rawLabel.toLowerCase());\n  if (newLabel === undefined) continue;\n\n  // Find the sub-range covering just the label text, then split the run by\n  // replacing its text (title-cased) and bolding that sub-range only.\n  const hits = para.search(rawLabel, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length === 0) continue;\n\n  const labelRange = hits.items[0];\n  labelRange.insertText(newLabel, \"Replace\");\n  labelRange.font.bold = true;\n}\n\nawait context.sync();\n", "ps1": "# Make the label at the start of each \"Label: value\" paragraph bold, and\n# normalize the label to Title Case (commit: \"modify the paragraphs to make\n# the title bold on each section\").\n#\n# Each target paragraph currently holds a single run \"Label: value\". We use\n# Find to seat a Range on just the label text, replace its text with the\n# title-cased label, and bold that sub-range \u2014 leaving \": value\" untouched in\n# a separate, non-bold run.\n\n$d = $word.ActiveDocument\n\n$labels = @(\n    \"Project title\",\n    \"Project description\",\n    \"Business driver\",\n    \"Business value\",\n    \"Business risk\",\n    \"Budget expense\",\n    \"Internal hours\",\n    \"External hours\",\n    \"Solutions involvement\",\n    \"Solutions hours\",\n    \"PMO involvement\",\n    \"PMO hours\",\n    \"Total expected hours\"\n)\n\nfunction Get-TitleCase($label) {\n    $words = $label.Split(\" \")\n    $out = @()\n    foreach ($w in $words) {\n        if ($w.ToUpper() -eq \"PMO\") {\n            $out += \"PMO\"\n        } else {\n            $out += ($w.Substring(0,1).ToUpper() + $w.Substring(1).ToLower())\n        }\n    }\n    return [string]::Join(\" \", $out)\n}\n\n$labelMap = @{}\nforeach ($l in $labels) {\n    $labelMap[$l.ToLower()] = Get-TitleCase $l\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $text = $para.Range.Text\n    $sep = $text.IndexOf(\": \")\n    if ($sep -lt 0) { continue }\n\n    $rawLabel = $text.Substring(0, $sep)\n    $key = $rawLabel.ToLower()\n    if (-not $labelMap.ContainsKey($key)) { continue }\n    $newLabel = $labelMap[$key]\n\n    $r = $para.Range.Duplicate\n    $found = $r.Find.Execute($rawLabel, $true)\n    if (-not $found) { continue }\n\n    $r.Text = $newLabel\n    $r.Font.Bold = 1\n}\n"}
